$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 15.981972
$ws.Range("H2").Value = 47.945916
$ws.Range("I2").Value = 0.1372507760882863
$ws.Range("J2").Value = 0.1372507760882863
$ws.Range("M2").Value = 23.641894
$ws.Range("N2").Value = 70.92568199999999
$ws.Range("O2").Value = 0.1609466983245457
$ws.Range("P2").Value = 0.1609466983245456
$ws.Range("Q2").Value = 377.8440879349679
$ws.Range("R2").Value = 3400.596791414711
$ws.Range("S2").Value = 0.02209005925389117
$ws.Range("T2").Value = 0.02209005925389117

# Row 3
$ws.Range("G3").Value = 15.981972
$ws.Range("H3").Value = 47.945916
$ws.Range("I3").Value = 0.1372507760882863
$ws.Range("J3").Value = 0.1372507760882863
$ws.Range("O3").Value = 0.2271397161392734
$ws.Range("P3").Value = 0.2271397161392734
$ws.Range("Q3").Value = 533.2411274780559
$ws.Range("R3").Value = 4799.170147302504
$ws.Range("S3").Value = 0.03117510232058831
$ws.Range("T3").Value = 0.03117510232058831

# Row 4
$ws.Range("G4").Value = 15.981972
$ws.Range("H4").Value = 47.945916
$ws.Range("I4").Value = 0.1372507760882863
$ws.Range("J4").Value = 0.1372507760882863
$ws.Range("M4").Value = 40.78693933333333
$ws.Range("N4").Value = 122.360818
$ws.Range("O4").Value = 0.2776648613881589
$ws.Range("P4").Value = 0.2776648613881589
$ws.Range("Q4").Value = 651.8557223910319
$ws.Range("R4").Value = 5866.701501519287
$ws.Range("S4").Value = 0.03810971771797124
$ws.Range("T4").Value = 0.03810971771797123

# Row 5
$ws.Range("G5").Value = 15.981972
$ws.Range("H5").Value = 47.945916
$ws.Range("I5").Value = 0.1372507760882863
$ws.Range("J5").Value = 0.1372507760882863
$ws.Range("M5").Value = 10.35975466666667
$ws.Range("N5").Value = 31.079264
$ws.Range("O5").Value = 0.07052600392558668
$ws.Range("P5").Value = 0.07052600392558667
$ws.Range("Q5").Value = 165.569309009536
$ws.Range("R5").Value = 1490.123781085824
$ws.Range("S5").Value = 0.009679748773192295
$ws.Range("T5").Value = 0.009679748773192293

# Row 6
$ws.Range("G6").Value = 15.981972
$ws.Range("H6").Value = 47.945916
$ws.Range("I6").Value = 0.1372507760882863
$ws.Range("J6").Value = 0.1372507760882863
$ws.Range("M6").Value = 38.73894066666667
$ws.Range("N6").Value = 116.216822
$ws.Range("O6").Value = 0.2637227202224355
$ws.Range("P6").Value = 0.2637227202224354
$ws.Range("Q6").Value = 619.124665044328
$ws.Range("R6").Value = 5572.121985398952
$ws.Range("S6").Value = 0.03619614802264325
$ws.Range("T6").Value = 0.03619614802264324

# Row 7
$ws.Range("H7").Value = 69.213024
$ws.Range("I7").Value = 0.1981303529463737
$ws.Range("J7").Value = 0.1981303529463737
$ws.Range("M7").Value = 23.641894
$ws.Range("N7").Value = 70.92568199999999
$ws.Range("O7").Value = 0.1609466983245457
$ws.Range("P7").Value = 0.1609466983245456
$ws.Range("Q7").Value = 545.442325609152
$ws.Range("R7").Value = 4908.980930482368
$ws.Range("S7").Value = 0.03188842614459576
$ws.Range("T7").Value = 0.03188842614459576

# Row 8
$ws.Range("H8").Value = 69.213024
$ws.Range("I8").Value = 0.1981303529463737
$ws.Range("J8").Value = 0.1981303529463737
$ws.Range("O8").Value = 0.2271397161392734
$ws.Range("P8").Value = 0.2271397161392734
$ws.Range("Q8").Value = 769.767980945984
$ws.Range("R8").Value = 6927.911828513857
$ws.Range("S8").Value = 0.04500327212681337
$ws.Range("T8").Value = 0.04500327212681336

# Row 9
$ws.Range("H9").Value = 69.213024
$ws.Range("I9").Value = 0.1981303529463737
$ws.Range("J9").Value = 0.1981303529463737
$ws.Range("M9").Value = 40.78693933333333
$ws.Range("N9").Value = 122.360818
$ws.Range("O9").Value = 0.2776648613881589
$ws.Range("P9").Value = 0.2776648613881589
$ws.Range("Q9").Value = 940.995803654848
$ws.Range("R9").Value = 8468.962232893633
$ws.Range("S9").Value = 0.05501383698764185
$ws.Range("T9").Value = 0.05501383698764184

# Row 10
$ws.Range("H10").Value = 69.213024
$ws.Range("I10").Value = 0.1981303529463737
$ws.Range("J10").Value = 0.1981303529463737
$ws.Range("M10").Value = 10.35975466666667
$ws.Range("N10").Value = 31.079264
$ws.Range("O10").Value = 0.07052600392558668
$ws.Range("P10").Value = 0.07052600392558667
$ws.Range("Q10").Value = 239.009982792704
$ws.Range("R10").Value = 2151.089845134336
$ws.Range("S10").Value = 0.01397334204967383
$ws.Range("T10").Value = 0.01397334204967382

# Row 11
$ws.Range("H11").Value = 69.213024
$ws.Range("I11").Value = 0.1981303529463737
$ws.Range("J11").Value = 0.1981303529463737
$ws.Range("M11").Value = 38.73894066666667
$ws.Range("N11").Value = 116.216822
$ws.Range("O11").Value = 0.2637227202224355
$ws.Range("P11").Value = 0.2637227202224354
$ws.Range("Q11").Value = 893.7464100321922
$ws.Range("R11").Value = 8043.717690289729
$ws.Range("S11").Value = 0.0522514756376489
$ws.Range("T11").Value = 0.05225147563764887

# Row 12
$ws.Range("G12").Value = 40.09539033333333
$ws.Range("H12").Value = 120.286171
$ws.Range("I12").Value = 0.3443331924754199
$ws.Range("J12").Value = 0.3443331924754199
$ws.Range("M12").Value = 23.641894
$ws.Range("N12").Value = 70.92568199999999
$ws.Range("O12").Value = 0.1609466983245457
$ws.Range("P12").Value = 0.1609466983245456
$ws.Range("Q12").Value = 947.9309681492913
$ws.Range("R12").Value = 8531.378713343622
$ws.Range("S12").Value = 0.05541929045246912
$ws.Range("T12").Value = 0.05541929045246911

# Row 13
$ws.Range("G13").Value = 40.09539033333333
$ws.Range("H13").Value = 120.286171
$ws.Range("I13").Value = 0.3443331924754199
$ws.Range("J13").Value = 0.3443331924754199
$ws.Range("O13").Value = 0.2271397161392734
$ws.Range("P13").Value = 0.2271397161392734
$ws.Range("Q13").Value = 1337.789300845942
$ws.Range("R13").Value = 12040.10370761347
$ws.Range("S13").Value = 0.07821174359619666
$ws.Range("T13").Value = 0.07821174359619665

# Row 14
$ws.Range("G14").Value = 40.09539033333333
$ws.Range("H14").Value = 120.286171
$ws.Range("I14").Value = 0.3443331924754199
$ws.Range("J14").Value = 0.3443331924754199
$ws.Range("M14").Value = 40.78693933333333
$ws.Range("N14").Value = 122.360818
$ws.Range("O14").Value = 0.2776648613881589
$ws.Range("P14").Value = 0.2776648613881589
$ws.Range("Q14").Value = 1635.368253071986
$ws.Range("R14").Value = 14718.31427764788
$ws.Range("S14").Value = 0.09560922816002972
$ws.Range("T14").Value = 0.09560922816002969

# Row 15
$ws.Range("G15").Value = 40.09539033333333
$ws.Range("H15").Value = 120.286171
$ws.Range("I15").Value = 0.3443331924754199
$ws.Range("J15").Value = 0.3443331924754199
$ws.Range("M15").Value = 10.35975466666667
$ws.Range("N15").Value = 31.079264
$ws.Range("O15").Value = 0.07052600392558668
$ws.Range("P15").Value = 0.07052600392558667
$ws.Range("Q15").Value = 415.3784071175715
$ws.Range("R15").Value = 3738.405664058144
$ws.Range("S15").Value = 0.02428444408423126
$ws.Range("T15").Value = 0.02428444408423125

# Row 16
$ws.Range("G16").Value = 40.09539033333333
$ws.Range("H16").Value = 120.286171
$ws.Range("I16").Value = 0.3443331924754199
$ws.Range("J16").Value = 0.3443331924754199
$ws.Range("M16").Value = 38.73894066666667
$ws.Range("N16").Value = 116.216822
$ws.Range("O16").Value = 0.2637227202224355
$ws.Range("P16").Value = 0.2637227202224354
$ws.Range("Q16").Value = 1553.252947129841
$ws.Range("R16").Value = 13979.27652416856
$ws.Range("S16").Value = 0.09080848618249318
$ws.Range("T16").Value = 0.09080848618249314

# Row 17
$ws.Range("G17").Value = 8.831340666666666
$ws.Range("H17").Value = 26.494022
$ws.Range("I17").Value = 0.07584222775512579
$ws.Range("J17").Value = 0.07584222775512579
$ws.Range("M17").Value = 23.641894
$ws.Range("N17").Value = 70.92568199999999
$ws.Range("O17").Value = 0.1609466983245457
$ws.Range("P17").Value = 0.1609466983245456
$ws.Range("Q17").Value = 208.7896199192226
$ws.Range("R17").Value = 1879.106579273004
$ws.Range("S17").Value = 0.01220655615076572
$ws.Range("T17").Value = 0.01220655615076571

# Row 18
$ws.Range("G18").Value = 8.831340666666666
$ws.Range("H18").Value = 26.494022
$ws.Range("I18").Value = 0.07584222775512579
$ws.Range("J18").Value = 0.07584222775512579
$ws.Range("O18").Value = 0.2271397161392734
$ws.Range("P18").Value = 0.2271397161392734
$ws.Range("Q18").Value = 294.6591355707631
$ws.Range("R18").Value = 2651.932220136868
$ws.Range("S18").Value = 0.01722678208366939
$ws.Range("T18").Value = 0.01722678208366939

# Row 19
$ws.Range("G19").Value = 8.831340666666666
$ws.Range("H19").Value = 26.494022
$ws.Range("I19").Value = 0.07584222775512579
$ws.Range("J19").Value = 0.07584222775512579
$ws.Range("M19").Value = 40.78693933333333
$ws.Range("N19").Value = 122.360818
$ws.Range("O19").Value = 0.2776648613881589
$ws.Range("P19").Value = 0.2776648613881589
$ws.Range("Q19").Value = 360.2033560033328
$ws.Range("R19").Value = 3241.830204029995
$ws.Range("S19").Value = 0.02105872165699618
$ws.Range("T19").Value = 0.02105872165699618

# Row 20
$ws.Range("G20").Value = 8.831340666666666
$ws.Range("H20").Value = 26.494022
$ws.Range("I20").Value = 0.07584222775512579
$ws.Range("J20").Value = 0.07584222775512579
$ws.Range("M20").Value = 10.35975466666667
$ws.Range("N20").Value = 31.079264
$ws.Range("O20").Value = 0.07052600392558668
$ws.Range("P20").Value = 0.07052600392558667
$ws.Range("Q20").Value = 91.4905226844231
$ws.Range("R20").Value = 823.4147041598079
$ws.Range("S20").Value = 0.005348849252383241
$ws.Range("T20").Value = 0.005348849252383239

# Row 21
$ws.Range("G21").Value = 8.831340666666666
$ws.Range("H21").Value = 26.494022
$ws.Range("I21").Value = 0.07584222775512579
$ws.Range("J21").Value = 0.07584222775512579
$ws.Range("M21").Value = 38.73894066666667
$ws.Range("N21").Value = 116.216822
$ws.Range("O21").Value = 0.2637227202224355
$ws.Range("P21").Value = 0.2637227202224354
$ws.Range("Q21").Value = 342.1167820931204
$ws.Range("R21").Value = 3079.051038838084
$ws.Range("S21").Value = 0.02000131861131127
$ws.Range("T21").Value = 0.02000131861131126

# Row 22
$ws.Range("G22").Value = 28.463871
$ws.Range("H22").Value = 85.39161300000001
$ws.Range("I22").Value = 0.2444434507347945
$ws.Range("J22").Value = 0.2444434507347945
$ws.Range("M22").Value = 23.641894
$ws.Range("N22").Value = 70.92568199999999
$ws.Range("O22").Value = 0.1609466983245457
$ws.Range("P22").Value = 0.1609466983245456
$ws.Range("Q22").Value = 672.939821011674
$ws.Range("R22").Value = 6056.458389105066
$ws.Range("S22").Value = 0.03934236632282391
$ws.Range("T22").Value = 0.0393423663228239

# Row 23
$ws.Range("G23").Value = 28.463871
$ws.Range("H23").Value = 85.39161300000001
$ws.Range("I23").Value = 0.2444434507347945
$ws.Range("J23").Value = 0.2444434507347945
$ws.Range("O23").Value = 0.2271397161392734
$ws.Range("P23").Value = 0.2271397161392734
$ws.Range("Q23").Value = 949.7017429657579
$ws.Range("R23").Value = 8547.315686691823
$ws.Range("S23").Value = 0.05552281601200568
$ws.Range("T23").Value = 0.05552281601200568

# Row 24
$ws.Range("G24").Value = 28.463871
$ws.Range("H24").Value = 85.39161300000001
$ws.Range("I24").Value = 0.2444434507347945
$ws.Range("J24").Value = 0.2444434507347945
$ws.Range("M24").Value = 40.78693933333333
$ws.Range("N24").Value = 122.360818
$ws.Range("O24").Value = 0.2776648613881589
$ws.Range("P24").Value = 0.2776648613881589
$ws.Range("Q24").Value = 1160.954179668826
$ws.Range("R24").Value = 10448.58761701944
$ws.Range("S24").Value = 0.06787335686551997
$ws.Range("T24").Value = 0.06787335686551996

# Row 25
$ws.Range("G25").Value = 28.463871
$ws.Range("H25").Value = 85.39161300000001
$ws.Range("I25").Value = 0.2444434507347945
$ws.Range("J25").Value = 0.2444434507347945
$ws.Range("M25").Value = 10.35975466666667
$ws.Range("N25").Value = 31.079264
$ws.Range("O25").Value = 0.07052600392558668
$ws.Range("P25").Value = 0.07052600392558667
$ws.Range("Q25").Value = 294.878720423648
$ws.Range("R25").Value = 2653.908483812832
$ws.Range("S25").Value = 0.01723961976610607
$ws.Range("T25").Value = 0.01723961976610607

# Row 26
$ws.Range("G26").Value = 28.463871
$ws.Range("H26").Value = 85.39161300000001
$ws.Range("I26").Value = 0.2444434507347945
$ws.Range("J26").Value = 0.2444434507347945
$ws.Range("M26").Value = 38.73894066666667
$ws.Range("N26").Value = 116.216822
$ws.Range("O26").Value = 0.2637227202224355
$ws.Range("P26").Value = 0.2637227202224354
$ws.Range("Q26").Value = 1102.660209812654
$ws.Range("R26").Value = 9923.941888313888
$ws.Range("S26").Value = 0.0644652917683389
$ws.Range("T26").Value = 0.06446529176833887
